$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newFilesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Boston Terrier']   
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

# Update the FilesTab query cell (B4) to remove the File Type and Breed columns
$ws.Cells.Item(4, 2).Value = $newFilesQuery

# The removed lines shrink the wrapped text, so the row shrinks to fit it
$ws.Rows.Item(4).RowHeight = 217.5

# Update the view/selection state to match the post-edit UI state
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("B4").Select()
